$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.002166666666666667
$ws.Range("H2").Value = 0.0065
$ws.Range("I2").Value = 0.004890446475191893
$ws.Range("J2").Value = 0.004890446475191893
$ws.Range("M2").Value = 0.6415476666666667
$ws.Range("N2").Value = 1.924643
$ws.Range("O2").Value = 0.1426849042655057
$ws.Range("P2").Value = 0.1426849042655057
$ws.Range("Q2").Value = 0.001390019944444444
$ws.Range("R2").Value = 0.0125101795
$ws.Range("S2").Value = 0.0006977928871283352
$ws.Range("T2").Value = 0.0006977928871283352

# Row 3
$ws.Range("G3").Value = 0.002166666666666667
$ws.Range("H3").Value = 0.0065
$ws.Range("I3").Value = 0.004890446475191893
$ws.Range("J3").Value = 0.004890446475191893
$ws.Range("O3").Value = 0.5986102210699216
$ws.Range("P3").Value = 0.5986102210699217
$ws.Range("Q3").Value = 0.005831592
$ws.Range("R3").Value = 0.052484328
$ws.Range("S3").Value = 0.002927471245645238
$ws.Range("T3").Value = 0.002927471245645239

# Row 4
$ws.Range("G4").Value = 0.002166666666666667
$ws.Range("H4").Value = 0.0065
$ws.Range("I4").Value = 0.004890446475191893
$ws.Range("J4").Value = 0.004890446475191893
$ws.Range("M4").Value = 1.163203
$ws.Range("N4").Value = 3.489609
$ws.Range("O4").Value = 0.2587048746645726
$ws.Range("P4").Value = 0.2587048746645726
$ws.Range("Q4").Value = 0.002520273166666667
$ws.Range("R4").Value = 0.0226824585
$ws.Range("S4").Value = 0.00126518234241832
$ws.Range("T4").Value = 0.00126518234241832

# Row 5
$ws.Range("G5").Value = 0.440874
$ws.Range("H5").Value = 1.322622
$ws.Range("I5").Value = 0.9951095535248081
$ws.Range("J5").Value = 0.9951095535248081
$ws.Range("M5").Value = 0.6415476666666667
$ws.Range("N5").Value = 1.924643
$ws.Range("O5").Value = 0.1426849042655057
$ws.Range("P5").Value = 0.1426849042655057
$ws.Range("Q5").Value = 0.282841685994
$ws.Range("R5").Value = 2.545575173946
$ws.Range("S5").Value = 0.1419871113783774
$ws.Range("T5").Value = 0.1419871113783774

# Row 6
$ws.Range("G6").Value = 0.440874
$ws.Range("H6").Value = 1.322622
$ws.Range("I6").Value = 0.9951095535248081
$ws.Range("J6").Value = 0.9951095535248081
$ws.Range("O6").Value = 0.5986102210699216
$ws.Range("P6").Value = 0.5986102210699217
$ws.Range("Q6").Value = 1.186614134496
$ws.Range("R6").Value = 10.679527210464
$ws.Range("S6").Value = 0.5956827498242764
$ws.Range("T6").Value = 0.5956827498242765

# Row 7
$ws.Range("G7").Value = 0.440874
$ws.Range("H7").Value = 1.322622
$ws.Range("I7").Value = 0.9951095535248081
$ws.Range("J7").Value = 0.9951095535248081
$ws.Range("M7").Value = 1.163203
$ws.Range("N7").Value = 3.489609
$ws.Range("O7").Value = 0.2587048746645726
$ws.Range("P7").Value = 0.2587048746645726
$ws.Range("Q7").Value = 0.5128259594219999
$ws.Range("R7").Value = 4.615433634797999
$ws.Range("S7").Value = 0.2574396923221544
$ws.Range("T7").Value = 0.2574396923221544
